$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Data5"
$ws.Range("F2").Value = "maciek@example.com"
$ws.Range("F2").Interior.Pattern = -4142

$ws.Range("F9").Select()
